# Included parameters & randomized order
# The header label in C1 ("image_title_1") is replaced with "image_title_B".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C1").Value = "image_title_B"

# Move/collapse the selection to E16, matching the saved view state.
$ws.Range("E16").Select()
